$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns
$ws.Range("A42").Value = 111908910
$ws.Range("B42").Value = 56969
$ws.Range("E42").Value = 102120
$ws.Range("Q42").Value = 485046
$ws.Range("R42").Value = 6653165
$ws.Range("S42").Value = 37

# Plain text columns
$ws.Range("C42").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("D42").Value = "VU"
$ws.Range("F42").Value = "Rödstrupig piplärka"
$ws.Range("G42").Value = "Anthus cervinus"
$ws.Range("H42").Value = "(Pallas, 1811)"
$ws.Range("M42").Value = "rastande"
$ws.Range("P42").Value = "Lönnfallet, Grängesberg, Dlr"
$ws.Range("T42").Value = "Dalarna"
$ws.Range("U42").Value = "Ludvika"
$ws.Range("V42").Value = "Dalarna"
$ws.Range("W42").Value = "Grangärde"
$ws.Range("AW42").Value = "Urban Grenmyr"
$ws.Range("AX42").Value = "Sten-Erik Bohlin"

# Text-like values that Excel would otherwise auto-convert (numbers/dates) -
# force them to be stored as literal text, matching the source workbook.
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = "2"
$ws.Range("I42").ClearFormats()

$ws.Range("Y42").NumberFormat = "@"
$ws.Range("Y42").Value = "2008-08-12"
$ws.Range("Y42").ClearFormats()

$ws.Range("AA42").NumberFormat = "@"
$ws.Range("AA42").Value = "2008-08-12"
$ws.Range("AA42").ClearFormats()

# Boolean columns
$ws.Range("AD42").Value = $false
$ws.Range("AE42").Value = $false
$ws.Range("AG42").Value = $false

# Empty-but-present text cells
$ws.Range("K42").Value = ""
$ws.Range("L42").Value = ""
$ws.Range("N42").Value = ""
$ws.Range("AT42").Value = ""
$ws.Range("AY42").Value = ""
